# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers table updates
$ws.Range("C4").Value = 106
$ws.Range("D4").Value = 98.8
$ws.Range("C5").Value = 110

# Good Drivers table: set Driver Vintage for the 22.150.3.1 row.
# A plain string assignment of "2022-08-29" gets auto-converted by Excel into a
# date serial number, so instead enter it as a text formula and then convert the
# formula to its literal value in place (Paste Values), which keeps the original
# cell style/number format intact and leaves behind a plain text value.
$cell = $ws.Range("E13")
$cell.Formula = '="2022-08-29"'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false
